# Fix date strings in column D (ExamDate) that were stored in MM-DD-YYYY
# format - convert them to DD-MM-YYYY format. These values are stored as
# plain text, so we must force a text write (via an apostrophe prefix on
# Value2) to stop Excel from re-parsing the digits as a real date.
#
# Mapping (old -> new):
#   04-10-2024 -> 10-04-2024
#   04-12-2024 -> 12-04-2024
#   04-09-2024 -> 09-04-2024
#   04-08-2024 -> 08-04-2024

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "04-10-2024" = "10-04-2024"
    "04-12-2024" = "12-04-2024"
    "04-09-2024" = "09-04-2024"
    "04-08-2024" = "08-04-2024"
}

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($map.ContainsKey($val)) {
        $newVal = $map[$val]
        $cell.Value2 = "'" + $newVal
        $cell.Style = "Normal"
    }
}
